$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove column D ("Effective end date") entirely. This shifts the old
# column E ("Effective End Time") left into D, and drops the now-unused
# "Effective end date" shared string. Clearing formats first avoids
# leaving a stale/orphaned column-width record behind after the delete.
$ws.Columns.Item(4).ClearFormats()
$ws.Columns.Item(4).Delete()

# Update the active selection to D7 to match the target workbook state
$ws.Range("D7").Select()
